$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (price + 1h volume change) scraped refresh.
# Each entry is the target cell and its new text value. A handful of "Price"
# cells are plain decimals (e.g. "7.61") that Excel would otherwise auto-detect
# and coerce into a Number; those are written with a leading apostrophe to force
# text entry and then have their style reset to Normal so no formatting residue
# (quote-prefix style) is left behind, matching the original plain-text cells.
$updates = @(
    @{ Cell = 'D2'; Value = '63.597.84'; ForceText = $false }
    @{ Cell = 'E2'; Value = '  -1.08%  '; ForceText = $false }
    @{ Cell = 'D3'; Value = '3.417.25'; ForceText = $false }
    @{ Cell = 'E3'; Value = '  -2.28%  '; ForceText = $false }
    @{ Cell = 'E4'; Value = '  +0.08%  '; ForceText = $false }
    @{ Cell = 'D5'; Value = '581.03'; ForceText = $true }
    @{ Cell = 'E5'; Value = '  -1.28%  '; ForceText = $false }
    @{ Cell = 'D6'; Value = '129.72'; ForceText = $true }
    @{ Cell = 'E6'; Value = '  -3.33%  '; ForceText = $false }
    @{ Cell = 'E7'; Value = '  +0.06%  '; ForceText = $false }
    @{ Cell = 'E8'; Value = '  -1.53%  '; ForceText = $false }
    @{ Cell = 'D9'; Value = '7.61'; ForceText = $true }
    @{ Cell = 'E9'; Value = '  +4.47%  '; ForceText = $false }
    @{ Cell = 'D10'; Value = '0.125'; ForceText = $true }
    @{ Cell = 'E10'; Value = '  +0.35%  '; ForceText = $false }
    @{ Cell = 'D11'; Value = '0.383'; ForceText = $true }
    @{ Cell = 'E11'; Value = '  -0.85%  '; ForceText = $false }
    @{ Cell = 'D12'; Value = '4.016.15'; ForceText = $false }
    @{ Cell = 'E12'; Value = '  -1.84%  '; ForceText = $false }
    @{ Cell = 'E13'; Value = '  -0.38%  '; ForceText = $false }
    @{ Cell = 'E14'; Value = '  -2.26%  '; ForceText = $false }
    @{ Cell = 'D15'; Value = '3.429.82'; ForceText = $false }
    @{ Cell = 'E15'; Value = '  -1.90%  '; ForceText = $false }
    @{ Cell = 'D16'; Value = '63.558.20'; ForceText = $false }
    @{ Cell = 'E16'; Value = '  -1.26%  '; ForceText = $false }
    @{ Cell = 'D17'; Value = '24.93'; ForceText = $true }
    @{ Cell = 'E17'; Value = '  -3.09%  '; ForceText = $false }
    @{ Cell = 'D18'; Value = '9.83'; ForceText = $true }
    @{ Cell = 'E18'; Value = '  -0.26%  '; ForceText = $false }
    @{ Cell = 'D19'; Value = '5.65'; ForceText = $true }
    @{ Cell = 'E19'; Value = '  -1.67%  '; ForceText = $false }
    @{ Cell = 'D20'; Value = '13.33'; ForceText = $true }
    @{ Cell = 'E20'; Value = '  -1.65%  '; ForceText = $false }
    @{ Cell = 'D21'; Value = '385.60'; ForceText = $true }
    @{ Cell = 'E21'; Value = '  -2.05%  '; ForceText = $false }
    @{ Cell = 'D22'; Value = '0.562'; ForceText = $true }
    @{ Cell = 'E22'; Value = '  -1.74%  '; ForceText = $false }
    @{ Cell = 'D23'; Value = '3.566.15'; ForceText = $false }
    @{ Cell = 'E23'; Value = '  -1.93%  '; ForceText = $false }
    @{ Cell = 'D24'; Value = '73.64'; ForceText = $true }
    @{ Cell = 'E24'; Value = '  -1.37%  '; ForceText = $false }
    @{ Cell = 'D25'; Value = '0.995'; ForceText = $true }
    @{ Cell = 'E25'; Value = '  -0.47%  '; ForceText = $false }
    @{ Cell = 'D26'; Value = '5.52'; ForceText = $true }
    @{ Cell = 'E26'; Value = '  -3.82%  '; ForceText = $false }
    @{ Cell = 'D27'; Value = '0.0000110'; ForceText = $true }
    @{ Cell = 'E27'; Value = '  -5.05%  '; ForceText = $false }
    @{ Cell = 'D28'; Value = '0.991'; ForceText = $true }
    @{ Cell = 'E28'; Value = '  -0.73%  '; ForceText = $false }
    @{ Cell = 'E29'; Value = '  -1.43%  '; ForceText = $false }
    @{ Cell = 'D30'; Value = '7.02'; ForceText = $true }
    @{ Cell = 'E30'; Value = '  -4.72%  '; ForceText = $false }
    @{ Cell = 'B31'; Value = 'InternetComputer(DFINITY)'; ForceText = $false }
    @{ Cell = 'C31'; Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; ForceText = $false }
    @{ Cell = 'D31'; Value = '7.94'; ForceText = $true }
    @{ Cell = 'E31'; Value = '  -3.63%  '; ForceText = $false }
    @{ Cell = 'B32'; Value = 'Kaspa'; ForceText = $false }
    @{ Cell = 'C32'; Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'; ForceText = $false }
    @{ Cell = 'D32'; Value = '0.155'; ForceText = $true }
    @{ Cell = 'E32'; Value = '  +2.19%  '; ForceText = $false }
    @{ Cell = 'E33'; Value = '  -4.08%  '; ForceText = $false }
    @{ Cell = 'D34'; Value = '3.455.43'; ForceText = $false }
    @{ Cell = 'E34'; Value = '  -1.79%  '; ForceText = $false }
    @{ Cell = 'D36'; Value = '22.83'; ForceText = $true }
    @{ Cell = 'E36'; Value = '  -2.65%  '; ForceText = $false }
    @{ Cell = 'D37'; Value = '5.17'; ForceText = $true }
    @{ Cell = 'E37'; Value = '  +0.38%  '; ForceText = $false }
    @{ Cell = 'D38'; Value = '6.76'; ForceText = $true }
    @{ Cell = 'E38'; Value = '  -1.76%  '; ForceText = $false }
    @{ Cell = 'D39'; Value = '163.63'; ForceText = $true }
    @{ Cell = 'E39'; Value = '  -2.26%  '; ForceText = $false }
    @{ Cell = 'E40'; Value = '  -3.12%  '; ForceText = $false }
    @{ Cell = 'D41'; Value = '0.0774'; ForceText = $true }
    @{ Cell = 'E41'; Value = '  -0.65%  '; ForceText = $false }
    @{ Cell = 'D42'; Value = '0.784'; ForceText = $true }
    @{ Cell = 'E42'; Value = '  -2.85%  '; ForceText = $false }
    @{ Cell = 'E43'; Value = '  +0.11%  '; ForceText = $false }
    @{ Cell = 'D44'; Value = '41.30'; ForceText = $true }
    @{ Cell = 'E44'; Value = '  -1.45%  '; ForceText = $false }
    @{ Cell = 'D45'; Value = '4.33'; ForceText = $true }
    @{ Cell = 'E45'; Value = '  -1.48%  '; ForceText = $false }
    @{ Cell = 'E46'; Value = '  -2.49%  '; ForceText = $false }
    @{ Cell = 'D47'; Value = '23.40'; ForceText = $true }
    @{ Cell = 'E47'; Value = '  -7.91%  '; ForceText = $false }
    @{ Cell = 'E48'; Value = '  -4.73%  '; ForceText = $false }
    @{ Cell = 'D49'; Value = '6.71'; ForceText = $true }
    @{ Cell = 'E49'; Value = '  -0.50%  '; ForceText = $false }
    @{ Cell = 'D50'; Value = '0.893'; ForceText = $true }
    @{ Cell = 'E50'; Value = '  +0.15%  '; ForceText = $false }
    @{ Cell = 'D51'; Value = '2.306.15'; ForceText = $false }
    @{ Cell = 'E51'; Value = '  -6.64%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $r = $ws.Range($u.Cell)
    if ($u.ForceText) {
        $r.Value = "'" + $u.Value
        $r.Style = "Normal"
    } else {
        $r.Value = $u.Value
    }
}

